$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q3" sheet right after "总计" (position 1) ---
# Clone the existing "2022-Q1" sheet (same column layout / styles) and
# drop it into the new slot, then overwrite its data with the 2022-Q3 figures.
$firstSheet = $wb.Worksheets.Item(1)
$templateSheet = $wb.Worksheets.Item("2022-Q1")
$templateSheet.Copy($null, $firstSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Columns B:G hold text values in this workbook (fund code/name/ratios as
# strings) - force text format so numeric-looking strings (e.g. "008115",
# "2.44") keep their literal text representation instead of being coerced
# into numbers.
$q3.Range("B2:G4").NumberFormat = "@"

$q3.Cells.Item(2,2).Value = "008115"
$q3.Cells.Item(2,3).Value = "天弘中证红利低波动100指数C"
$q3.Cells.Item(2,4).Value = "2.44"
$q3.Cells.Item(2,5).Value = "94.56"
$q3.Cells.Item(2,6).Value = "1.81"
$q3.Cells.Item(2,7).Value = "0.0442"
$q3.Cells.Item(2,8).Value = 6

$q3.Cells.Item(3,2).Value = "008114"
$q3.Cells.Item(3,3).Value = "天弘中证红利低波动100指数A"
$q3.Cells.Item(3,4).Value = "1.89"
$q3.Cells.Item(3,5).Value = "94.56"
$q3.Cells.Item(3,6).Value = "1.81"
$q3.Cells.Item(3,7).Value = "0.0342"
$q3.Cells.Item(3,8).Value = 6

$q3.Cells.Item(4,2).Value = "515100"
$q3.Cells.Item(4,3).Value = "景顺长城中证红利低波动100ETF"
$q3.Cells.Item(4,4).Value = "1.62"
$q3.Cells.Item(4,5).Value = "98.63"
$q3.Cells.Item(4,6).Value = "1.90"
$q3.Cells.Item(4,7).Value = "0.0308"
$q3.Cells.Item(4,8).Value = 6

# --- 2. Update the "总计" summary sheet: insert a row for 2022-Q3 and
#        shift the existing quarters down by one ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
# Re-apply the data-row formatting (lost by the blank row Insert() leaves
# behind) by copying it from the row below (still holding the old row-2
# content/format at this point).
$summary.Range("A3:D3").Copy($summary.Range("A2:D2"))

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.11

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q1"
$summary.Cells.Item(3,3).Value = 3
$summary.Cells.Item(3,4).Value = 0.45

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2021-Q4"
$summary.Cells.Item(4,3).Value = 3
$summary.Cells.Item(4,4).Value = 0.45

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q3"
$summary.Cells.Item(5,3).Value = 8
$summary.Cells.Item(5,4).Value = 4.03

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q2"
$summary.Cells.Item(6,3).Value = 17
$summary.Cells.Item(6,4).Value = 7.89

$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q1"
$summary.Cells.Item(7,3).Value = 15
$summary.Cells.Item(7,4).Value = 7.66

$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(8,2).Value = "2020-Q4"
$summary.Cells.Item(8,3).Value = 22
$summary.Cells.Item(8,4).Value = 10.6

# --- 3. Restore the originally active sheet ("2020-Q4" was the selected
#        tab before this edit; our sheet-copy/rename above made the new
#        "2022-Q3" sheet active, so switch back) ---
$wb.Worksheets.Item("2020-Q4").Activate()
